# Card Chronicles: Prophecy - 卡牌.xlsx
# Add an "英文名" (English name) column after column A, populate it with the
# English codename for each card, and drop the old "是否完成" (done?) column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Structural change: insert a new column B for the English names, and
#    remove the old "是否完成" / "是" column (originally H, now shifted to I).
# ---------------------------------------------------------------------------
$ws.Columns("B").Insert()
$ws.Columns("I").Delete()

# ---------------------------------------------------------------------------
# 2. Header + English-name content for column B
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "英文名"

$englishNames = @{
    2  = "Rain"
    3  = "CalamityGone"
    4  = "fog"
    5  = "happy"
    6  = "genshin"
    7  = "song"
    8  = "untie"
    9  = "revival"
    10 = "Doppelganger"
    11 = "armour"
    12 = "Sacrifice"
    13 = "beat"
    14 = "circlebeat"
    15 = "disablebeat"
    16 = "thunder"
    17 = "heavybeat"
    18 = "Sacrificebeat"
    19 = "fire"
    20 = "duantoutai"
    21 = "summon"
    22 = "paw"
    23 = "shenwei"
    24 = "poison"
}

foreach ($row in $englishNames.Keys) {
    $ws.Range("B$row").Value = $englishNames[$row]
}

# ---------------------------------------------------------------------------
# 3. Row heights that visibly changed in the edit (content/wrap driven)
# ---------------------------------------------------------------------------
$rowHeights = @{
    3  = 90
    4  = 33
    5  = 33
    6  = 60
    11 = 33
    17 = 33
    18 = 60
    21 = 66
    22 = 66
    23 = 33
    25 = 33
    26 = 33
    27 = 33
    28 = 33
    29 = 33
    30 = 33
    31 = 33
}

foreach ($row in $rowHeights.Keys) {
    $ws.Rows($row).RowHeight = $rowHeights[$row]
}

# ---------------------------------------------------------------------------
# 4. Column widths (new English-name column + slight re-tune of the rest)
# ---------------------------------------------------------------------------
$ws.Columns("B").ColumnWidth = 20.29
$ws.Columns("C").ColumnWidth = 17.71
$ws.Columns("D").ColumnWidth = 143.57
$ws.Columns("F").ColumnWidth = 66.14
$ws.Columns("G").ColumnWidth = 166.14
$ws.Columns("H").ColumnWidth = 16.14

# ---------------------------------------------------------------------------
# 5. View: unfreeze panes, zoom out, move the selection
# ---------------------------------------------------------------------------
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.Zoom = 70
$ws.Range("C25").Select()
